$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update the policy number on row 2 (NroPoliza); leading apostrophe keeps
# it stored as text (matches the existing Text number format on the cell)
$ws.Range("E2").Value = "'04104013566"

# Remove all hyperlinks on the sheet (B3 / B5 used to link out)
$ws.Hyperlinks.Delete()

# Clear out the old test rows (3 and 5), leaving only the formatted
# placeholder cells behind (B and E keep their style, rest are blank)
$ws.Range("A3").ClearContents()
$ws.Range("C3").ClearContents()
$ws.Range("D3").ClearContents()
$ws.Range("B3").ClearContents()
$ws.Range("E3").ClearContents()

$ws.Range("A5").ClearContents()
$ws.Range("C5").ClearContents()
$ws.Range("D5").ClearContents()
$ws.Range("B5").ClearContents()
$ws.Range("E5").ClearContents()

# Move the active selection to E2
$ws.Range("E2").Select()
